$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = -12.9717
$ws.Range("B8").Value = 5.959299999999998
$ws.Range("B10").Value = 5.566000000000003
$ws.Range("B12").Value = 5.0678
$ws.Range("C15").Value = -14.11419999999999
$ws.Range("B18").Value = 7.449199999999998
$ws.Range("C18").Value = -12.40569999999999
$ws.Range("C20").Value = -12.0002
$ws.Range("C29").Value = -11.3523
$ws.Range("C30").Value = -13.20039999999999
$ws.Range("C31").Value = -13.20279999999999
$ws.Range("B37").Value = 8.755699999999997
$ws.Range("C40").Value = -12.9864
$ws.Range("C50").Value = -13.53329999999999
$ws.Range("B55").Value = 6.271199999999997
$ws.Range("B68").Value = 6.131
$ws.Range("C68").Value = -11.94840000000001
$ws.Range("C76").Value = -12.28910000000001
$ws.Range("B77").Value = 9.039300000000003
$ws.Range("B78").Value = 9.4184
$ws.Range("B81").Value = 5.357199999999999
$ws.Range("B82").Value = 5.306900000000001
$ws.Range("C87").Value = -14.02659999999999
$ws.Range("C88").Value = -12.85659999999999
$ws.Range("C96").Value = -13.3606
$ws.Range("C98").Value = -12.61929999999999
$ws.Range("C101").Value = -12.88220000000001
$ws.Range("C102").Value = -13.2518
